$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" '43.229.31'
$ws.Range("E2").Value = '  -1.57%  '
Set-TextValue "D3" '2.360.61'
$ws.Range("E3").Value = '  +4.37%  '
$ws.Range("E4").Value = '  +0.08%  '
Set-TextValue "D5" '233.08'
$ws.Range("E5").Value = '  +0.23%  '
Set-TextValue "D6" '0.649'
$ws.Range("E6").Value = '  -1.15%  '
Set-TextValue "D7" '72.05'
$ws.Range("E7").Value = '  +13.16%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +6.52%  '
Set-TextValue "D10" '0.0986'
$ws.Range("E10").Value = '  +0.79%  '
Set-TextValue "D11" '27.14'
$ws.Range("E11").Value = '  +1.49%  '
Set-TextValue "D12" '2.717.20'
$ws.Range("E12").Value = '  +4.61%  '
$ws.Range("E13").Value = '  +0.86%  '
Set-TextValue "D14" '16.02'
$ws.Range("E14").Value = '  +2.47%  '
Set-TextValue "D15" '6.26'
$ws.Range("E15").Value = '  +1.99%  '
Set-TextValue "D16" '0.863'
$ws.Range("E16").Value = '  +2.56%  '
Set-TextValue "D17" '2.363.61'
$ws.Range("E17").Value = '  +4.45%  '
Set-TextValue "D18" '43.281.46'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("E19").Value = '  +3.67%  '
Set-TextValue "D20" '6.33'
$ws.Range("E20").Value = '  +3.09%  '
Set-TextValue "D21" '74.31'
Set-TextValue "D22" '249.77'
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  -1.04%  '
Set-TextValue "D25" '2.44'
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("E26").Value = '  +1.50%  '
Set-TextValue "D27" '9.98'
$ws.Range("E27").Value = '  +0.61%  '
Set-TextValue "D28" '22.37'
$ws.Range("E28").Value = '  +1.94%  '
Set-TextValue "D29" '173.75'
$ws.Range("E29").Value = '  -0.41%  '
$ws.Range("E30").Value = '  +5.09%  '
Set-TextValue "D31" '0.129'
$ws.Range("E31").Value = '  -5.34%  '
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("E34").Value = '  +1.16%  '
Set-TextValue "D35" '5.03'
$ws.Range("E35").Value = '  +1.54%  '
$ws.Range("E36").Value = '  +6.48%  '
$ws.Range("E37").Value = '  +2.06%  '
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("E41").Value = '  +2.61%  '
Set-TextValue "D42" '18.56'
$ws.Range("E42").Value = '  +7.91%  '
$ws.Range("E43").Value = '  +7.73%  '
Set-TextValue "D44" '100.02'
$ws.Range("E44").Value = '  +1.36%  '
Set-TextValue "D45" '4.50'
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("E46").Value = '  +1.44%  '
Set-TextValue "D47" '0.0952'
$ws.Range("E47").Value = '  -0.23%  '
Set-TextValue "D48" '1.444.30'
$ws.Range("E48").Value = '  -0.85%  '
Set-TextValue "D49" '2.589.66'
$ws.Range("E49").Value = '  +4.71%  '
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D50" '2.77'
$ws.Range("E50").Value = '  +0.49%  '
$ws.Range("B51").Value = 'TerraClassic'
$ws.Range("C51").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextValue "D51" '0.000202'
$ws.Range("E51").Value = '  -3.78%  '
